$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "28.412.42"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.11%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.828.45"
Set-TextValue $ws.Cells.Item(3, 5) "  +0.63%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9947"
Set-TextValue $ws.Cells.Item(4, 5) "  -0.59%  "
Set-TextValue $ws.Cells.Item(5, 4) "327.19"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.58%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.9918"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.65%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.4445"
Set-TextValue $ws.Cells.Item(7, 5) "  +0.15%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3790"
Set-TextValue $ws.Cells.Item(8, 5) "  -0.83%  "
Set-TextValue $ws.Cells.Item(9, 4) "45.49"
Set-TextValue $ws.Cells.Item(9, 5) "  +1.54%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.07771"
Set-TextValue $ws.Cells.Item(10, 5) "  +0.25%  "
Set-TextValue $ws.Cells.Item(11, 4) "1.140"
Set-TextValue $ws.Cells.Item(11, 5) "  -1.38%  "
Set-TextValue $ws.Cells.Item(12, 4) "22.28"
Set-TextValue $ws.Cells.Item(12, 5) "  -2.81%  "
Set-TextValue $ws.Cells.Item(13, 4) "0.9920"
Set-TextValue $ws.Cells.Item(13, 5) "  -0.75%  "
Set-TextValue $ws.Cells.Item(14, 4) "6.315"
Set-TextValue $ws.Cells.Item(14, 5) "  -0.86%  "
Set-TextValue $ws.Cells.Item(15, 4) "7.548"
Set-TextValue $ws.Cells.Item(15, 5) "  -0.94%  "
Set-TextValue $ws.Cells.Item(16, 4) "1.821.63"
Set-TextValue $ws.Cells.Item(16, 5) "  +0.41%  "
Set-TextValue $ws.Cells.Item(17, 4) "92.47"
Set-TextValue $ws.Cells.Item(17, 5) "  +13.02%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.00001085"
Set-TextValue $ws.Cells.Item(18, 5) "  -1.32%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.06377"
Set-TextValue $ws.Cells.Item(19, 5) "  -5.56%  "
Set-TextValue $ws.Cells.Item(20, 4) "0.9934"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.49%  "
Set-TextValue $ws.Cells.Item(21, 4) "17.59"
Set-TextValue $ws.Cells.Item(21, 5) "  -1.72%  "
Set-TextValue $ws.Cells.Item(22, 4) "6.348"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.30%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.5364"
Set-TextValue $ws.Cells.Item(23, 5) "  -1.32%  "
Set-TextValue $ws.Cells.Item(24, 4) "28.444.26"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.94%  "
Set-TextValue $ws.Cells.Item(25, 4) "11.77"
Set-TextValue $ws.Cells.Item(25, 5) "  -1.51%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.170"
Set-TextValue $ws.Cells.Item(26, 5) "  -10.92%  "
Set-TextValue $ws.Cells.Item(27, 4) "21.01"
Set-TextValue $ws.Cells.Item(27, 5) "  +0.79%  "
Set-TextValue $ws.Cells.Item(28, 4) "153.80"
Set-TextValue $ws.Cells.Item(28, 5) "  +0.23%  "
Set-TextValue $ws.Cells.Item(29, 4) "2.380"
Set-TextValue $ws.Cells.Item(29, 5) "  -1.12%  "
Set-TextValue $ws.Cells.Item(30, 4) "2.025.50"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.36%  "
Set-TextValue $ws.Cells.Item(31, 4) "130.04"
Set-TextValue $ws.Cells.Item(31, 5) "  -2.55%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.220"
Set-TextValue $ws.Cells.Item(32, 5) "  -5.86%  "
Set-TextValue $ws.Cells.Item(33, 4) "5.875"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.62%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.09254"
Set-TextValue $ws.Cells.Item(34, 5) "  -0.89%  "
Set-TextValue $ws.Cells.Item(35, 4) "3.658"
Set-TextValue $ws.Cells.Item(35, 5) "  -7.91%  "
Set-TextValue $ws.Cells.Item(36, 4) "12.90"
Set-TextValue $ws.Cells.Item(36, 5) "  +4.05%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.02356"
Set-TextValue $ws.Cells.Item(37, 5) "  -0.20%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.2199"
Set-TextValue $ws.Cells.Item(38, 5) "  -3.86%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.6630"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.52%  "
Set-TextValue $ws.Cells.Item(40, 4) "5.188"
Set-TextValue $ws.Cells.Item(40, 5) "  -1.83%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.06237"
Set-TextValue $ws.Cells.Item(41, 5) "  -2.89%  "
Set-TextValue $ws.Cells.Item(42, 4) "1.191"
Set-TextValue $ws.Cells.Item(42, 5) "  -1.96%  "
Set-TextValue $ws.Cells.Item(43, 4) "8.064"
Set-TextValue $ws.Cells.Item(43, 5) "  -2.16%  "
Set-TextValue $ws.Cells.Item(44, 4) "1.410"
Set-TextValue $ws.Cells.Item(44, 5) "  -2.93%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.9919"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.59%  "
Set-TextValue $ws.Cells.Item(46, 4) "13.87"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.73%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.6125"
Set-TextValue $ws.Cells.Item(47, 5) "  -0.64%  "
Set-TextValue $ws.Cells.Item(48, 4) "3.747"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.76%  "
Set-TextValue $ws.Cells.Item(49, 4) "127.46"
Set-TextValue $ws.Cells.Item(49, 5) "  -1.73%  "
Set-TextValue $ws.Cells.Item(50, 4) "2.038"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.48%  "

# Row 51: Cronos -> Aave
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 3).Style = "Normal"
Set-TextValue $ws.Cells.Item(51, 4) "79.54"
Set-TextValue $ws.Cells.Item(51, 5) "  +0.47%  "
